$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new data
$ws.Range("A2").Value = "Ahmed Taoufiq"
$ws.Range("B2").Value = "BB196497"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "114379134964978564311346"
$ws.Range("D2").Value = "CIH TEST"
$ws.Range("E2").Value = "CIH"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "001/LF/DR01"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 46000

# Delete rows 3 to 7 (old leftover data)
$ws.Range("A3:K7").EntireRow.Delete()
